$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3343.027
$ws.Range("I70").Value = 759.9
$ws.Range("J70").Value = 4299.7407
$ws.Range("K70").Value = 2279.7
$ws.Range("L70").Value = 12899.2221
$ws.Range("M70").Value = -2009.7
$ws.Range("N70").Value = -13439.2221
$ws.Range("H73").Value = 3343.027
$ws.Range("I73").Value = 759.9
$ws.Range("J73").Value = 4299.7407
$ws.Range("K73").Value = 2279.7
$ws.Range("L73").Value = 12899.2221
$ws.Range("M73").Value = -1343.7
$ws.Range("N73").Value = -14771.2221
$ws.Range("H98").Value = 3749.2092
$ws.Range("I98").Value = 3511.9375
$ws.Range("K98").Value = 3511.9375
$ws.Range("M98").Value = -2013.9375
$ws.Range("H112").Value = 7390.2
$ws.Range("J112").Value = 7699.5527
$ws.Range("L112").Value = 23098.6581
$ws.Range("N112").Value = -25314.6581
$ws.Range("H118").Value = 750.46155
$ws.Range("I118").Value = 614.7
$ws.Range("J118").Value = 1203
$ws.Range("K118").Value = 1844.1
$ws.Range("L118").Value = 3609
$ws.Range("M118").Value = -187.1000000000001
$ws.Range("N118").Value = -6923
$ws.Range("H122").Value = 3749.2092
$ws.Range("I122").Value = 3511.9375
$ws.Range("K122").Value = 10535.8125
$ws.Range("M122").Value = -8085.8125
$ws.Range("H125").Value = 29461.727
$ws.Range("I125").Value = 51415.5
$ws.Range("K125").Value = 462739.5
$ws.Range("M125").Value = -460279.5
$ws.Range("H135").Value = 23434.615
$ws.Range("I135").Value = 1203.0476
$ws.Range("K135").Value = 10827.4284
$ws.Range("M135").Value = -8292.428400000001
$ws.Range("H137").Value = 13790.385
$ws.Range("I137").Value = 18441.223
$ws.Range("J137").Value = 3326
$ws.Range("K137").Value = 55323.66900000001
$ws.Range("L137").Value = 9978
$ws.Range("M137").Value = -52773.66900000001
$ws.Range("N137").Value = -15078
$ws.Range("H138").Value = 3097.4849
$ws.Range("I138").Value = 2286
$ws.Range("K138").Value = 6858
$ws.Range("M138").Value = -1718
$ws.Range("H141").Value = 3787.6365
$ws.Range("I141").Value = 3607.5881
$ws.Range("J141").Value = 4399.8
$ws.Range("K141").Value = 10822.7643
$ws.Range("L141").Value = 13199.4
$ws.Range("M141").Value = -5642.764299999999
$ws.Range("N141").Value = -23559.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3416.6667
$ws.Range("I61").Value = 1356.3889
$ws.Range("J61").Value = 5476.9443
$ws.Range("K61").Value = 1356.3889
$ws.Range("L61").Value = 5476.9443
$ws.Range("M61").Value = -1144.3889
$ws.Range("N61").Value = -5900.9443
$ws.Range("H74").Value = 287916.53
$ws.Range("I74").Value = 335069.34
$ws.Range("K74").Value = 335069.34
$ws.Range("M74").Value = -334195.34
$ws.Range("H77").Value = 287916.53
$ws.Range("I77").Value = 335069.34
$ws.Range("K77").Value = 1675346.7
$ws.Range("M77").Value = -1670978.7
$ws.Range("H97").Value = 1718
$ws.Range("J97").Value = 1979.5
$ws.Range("L97").Value = 1979.5
$ws.Range("N97").Value = -2971.5
$ws.Range("H102").Value = 5453.9736
$ws.Range("I102").Value = 5527.303
$ws.Range("K102").Value = 5527.303
$ws.Range("M102").Value = -3905.303
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 1564.4878
$ws.Range("I132").Value = 1430.6428
$ws.Range("J132").Value = 1852.7693
$ws.Range("K132").Value = 4291.928400000001
$ws.Range("L132").Value = 5558.3079
$ws.Range("M132").Value = -1761.928400000001
$ws.Range("N132").Value = -10618.3079
$ws.Range("H136").Value = 3416.6667
$ws.Range("I136").Value = 1356.3889
$ws.Range("J136").Value = 5476.9443
$ws.Range("K136").Value = 4069.1667
$ws.Range("L136").Value = 16430.8329
$ws.Range("M136").Value = -1519.1667
$ws.Range("N136").Value = -21530.8329

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 25519.857
$ws.Range("I107").Value = 31183.824
$ws.Range("K107").Value = 31183.824
$ws.Range("M107").Value = -29263.824

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7145111
$ws.Range("I31").Value = 9092464
$ws.Range("K31").Value = 9092464
$ws.Range("M31").Value = -9092169
$ws.Range("H34").Value = 7145111
$ws.Range("I34").Value = 9092464
$ws.Range("K34").Value = 9092464
$ws.Range("M34").Value = -9092262
$ws.Range("H59").Value = 40000
$ws.Range("J59").Value = 40000
$ws.Range("L59").Value = 40000
$ws.Range("N59").Value = -42290
$ws.Range("H74").Value = 83332
$ws.Range("J74").Value = 83332
$ws.Range("L74").Value = 83332
$ws.Range("N74").Value = -85080
$ws.Range("H76").Value = 111118490
$ws.Range("I76").Value = 111118490
$ws.Range("K76").Value = 111118490
$ws.Range("M76").Value = -111118175
$ws.Range("H77").Value = 83332
$ws.Range("J77").Value = 83332
$ws.Range("L77").Value = 249996
$ws.Range("N77").Value = -258732
$ws.Range("H79").Value = 111118490
$ws.Range("I79").Value = 111118490
$ws.Range("K79").Value = 111118490
$ws.Range("M79").Value = -111117398
$ws.Range("H99").Value = 7078.647
$ws.Range("I99").Value = 6747.6665
$ws.Range("K99").Value = 6747.6665
$ws.Range("M99").Value = -5249.6665
$ws.Range("H107").Value = 1618.6
$ws.Range("J107").Value = 3484.625
$ws.Range("L107").Value = 3484.625
$ws.Range("N107").Value = -7324.625
$ws.Range("H122").Value = 1907.3334
$ws.Range("I122").Value = 1907.3334
$ws.Range("K122").Value = 5722.0002
$ws.Range("M122").Value = -3272.0002
$ws.Range("H126").Value = 7078.647
$ws.Range("I126").Value = 6747.6665
$ws.Range("K126").Value = 20242.9995
$ws.Range("M126").Value = -17772.9995
$ws.Range("H134").Value = 2329.25
$ws.Range("I134").Value = 2284.2778
$ws.Range("K134").Value = 6852.8334
$ws.Range("M134").Value = -4317.8334
$ws.Range("H141").Value = 241706.58
$ws.Range("J141").Value = 241706.58
$ws.Range("L141").Value = 241706.58
$ws.Range("N141").Value = -252066.58

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3120
$ws.Range("I5").Value = 3120
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 9360
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -9248
$ws.Range("H23").Value = 296.875
$ws.Range("I23").Value = 93.75
$ws.Range("J23").Value = 500
$ws.Range("K23").Value = 281.25
$ws.Range("L23").Value = 1500
$ws.Range("M23").Value = -46.25
$ws.Range("N23").Value = -1970
$ws.Range("H68").Value = 3988.9312
$ws.Range("I68").Value = 899
$ws.Range("J68").Value = 4217.815
$ws.Range("K68").Value = 2697
$ws.Range("L68").Value = 12653.445
$ws.Range("M68").Value = -1886
$ws.Range("N68").Value = -14275.445
$ws.Range("H71").Value = 3988.9312
$ws.Range("I71").Value = 899
$ws.Range("J71").Value = 4217.815
$ws.Range("K71").Value = 8091
$ws.Range("L71").Value = 37960.335
$ws.Range("M71").Value = -4035
$ws.Range("N71").Value = -46072.335
$ws.Range("H135").Value = 3120
$ws.Range("I135").Value = 3120
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 28080
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -25545
$ws.Range("H136").Value = 2076.6667
$ws.Range("I136").Value = 1816.3334
$ws.Range("K136").Value = 5449.0002
$ws.Range("M136").Value = -349.0002000000004
$ws.Range("H140").Value = 2783.375
$ws.Range("I140").Value = 2252.7856
$ws.Range("K140").Value = 6758.3568
$ws.Range("M140").Value = -1578.3568

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 20000
$ws.Range("I132").Value = 20000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 60000
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -57470

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5337.9565
$ws.Range("I46").Value = 1886
$ws.Range("J46").Value = 6848.1875
$ws.Range("K46").Value = 1886
$ws.Range("L46").Value = 6848.1875
$ws.Range("M46").Value = -1698
$ws.Range("N46").Value = -7224.1875
$ws.Range("H100").Value = 2109.1333
$ws.Range("I100").Value = 2148.818
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 2148.818
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1607.818
$ws.Range("N100").Value = -3082
$ws.Range("H122").Value = 3772.1072
$ws.Range("I122").Value = 3985.158
$ws.Range("J122").Value = 3322.3333
$ws.Range("K122").Value = 11955.474
$ws.Range("L122").Value = 9966.999899999999
$ws.Range("M122").Value = -9505.474
$ws.Range("N122").Value = -14866.9999
$ws.Range("H136").Value = 4545.125
$ws.Range("I136").Value = 4489.0713
$ws.Range("K136").Value = 13467.2139
$ws.Range("M136").Value = -10917.2139

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 82040.66
$ws.Range("I122").Value = 106439.95
$ws.Range("J122").Value = 5357.143
$ws.Range("K122").Value = 319319.85
$ws.Range("L122").Value = 16071.429
$ws.Range("M122").Value = -316869.85
$ws.Range("N122").Value = -20971.429
$ws.Range("H132").Value = 6537.7607
$ws.Range("I132").Value = 6377.676
$ws.Range("K132").Value = 19133.028
$ws.Range("M132").Value = -16603.028
$ws.Range("H141").Value = 120403.17
$ws.Range("I141").Value = 115499
$ws.Range("J141").Value = 120870.234
$ws.Range("K141").Value = 115499
$ws.Range("L141").Value = 120870.234
$ws.Range("M141").Value = -110319
$ws.Range("N141").Value = -131230.234
